$d = $word.ActiveDocument

# The paragraph currently reads "Version 2." split across runs:
#   "Versi" | "on" | [spellEnd] | " 2" | [bookmark] | "."
# Target: "Version 1." split as:
#   "Version" | [spellEnd] | " 1." | [bookmark]

# 1) Merge "Versi" + "on" into a single "Version" run (chars 0-7).
#    Force an actual text change so the runs are rewritten/merged instead
#    of being left untouched as a no-op.
$d.Range(0, 7).Text = "Version "
$d.Range(0, 8).Text = "Version"

# 2) Remove the trailing "." run that follows the bookmark (chars 9-10).
$d.Range(9, 10).Text = ""

# 3) Change " 2" into " 1." (chars 7-9, the space+digit run).
$d.Range(7, 9).Text = " 1."
